$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128 (shifts old rows 128..183 down to 129..184)
$ws.Rows.Item(128).EntireRow.Insert()

# Populate the new row 128 with the new "Ají" record
$ws.Range("A128").Value = 7
$ws.Range("B128").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C128").Value = "Ñuble"
$ws.Range("D128").Value = 45027
$ws.Range("E128").Value = 16
$ws.Range("F128").Value = 100112021
$ws.Range("G128").Value = "Ají"
$ws.Range("H128").Value = "Cristal"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 60
$ws.Range("K128").Value = 15000
$ws.Range("L128").Value = 16000
$ws.Range("M128").Value = 15500
$ws.Range("N128").Value = "$/saco 25 kilos"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 620
$ws.Range("Q128").Value = 25
$ws.Range("R128").Value = "Hortaliza"
